# Edit slide 4 ("07-02a_BusinessRules4__closed_months.sql") - Content Placeholder 2
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

$CR = [char]13

# Work from the bottom of the affected block upward so paragraph indices
# for the parts we still need to touch stay valid.

# --- Paragraph 12 (COM index): "a product cannot change VAT percent ..." (lvl 1) -> delete
$tr.Paragraphs(12, 1).Delete()

# --- Paragraph 11: "invoices paid with the same receipt ..." (lvl 1) -> delete
$tr.Paragraphs(11, 1).Delete()

# --- Paragraph 10: "Other simple business rules" -> delete
$tr.Paragraphs(10, 1).Delete()

# --- Paragraph 9: "refusals.refusal_dt ..." (lvl 1, multi-run) -> delete, then
#     add a blank trailing paragraph (sz 2400, no indent level) after paragraph 8.
$tr.Paragraphs(9, 1).Delete()

# --- Paragraph 8: "receipt_details.receipt_date ..." (lvl 1, multi-run)
#     -> becomes single run "mutating tables (part II)<tab><tab>"
$para8 = $tr.Paragraphs(8, 1)
$para8.Text = "mutating tables (part II)`t`t"

# Insert the new trailing blank paragraph (ends the "Other technical stuff" list)
$blank2 = $para8.InsertAfter($CR)
$blank2.Font.Size = 24

# --- Paragraph 7: "cancelled_invoices.cancellation_dt ..." (lvl 1, multi-run)
#     -> becomes single run "at the end of a package body ... function <tab>"
$para7 = $tr.Paragraphs(7, 1)
$para7.Text = "at the end of a package body an anonymous block can be inserted; this will be executed once per session (at first call of a package variable, procedure or function `t"

# --- Paragraph 6: "Date/timestamp attribute synchronization for attributes in different tables " -> "Other technical stuff: "
$tr.Paragraphs(6, 1).Text = "Other technical stuff: "

# --- Paragraphs 3-5: former_contacts / former_vat_percents / receipts (lvl 1) -> delete entirely
$tr.Paragraphs(5, 1).Delete()
$tr.Paragraphs(4, 1).Delete()
$tr.Paragraphs(3, 1).Delete()

# --- Paragraph 2: "Date/timestamp attribute synchronization for attributes in the same table"
#     -> "Blocking/closing months (in Finance/Accounting)"
$para2 = $tr.Paragraphs(2, 1)
$para2.Text = "Blocking/closing months (in Finance/Accounting)"

# Insert the new "After the general ledger..." paragraph (3 runs, same sz/font)
$newPara = $para2.InsertAfter($CR + "After the general ledger, the VAT statements and other financial statements were finalized (that can happen monthly, quartely, twice a year or yearly), the operations incorporated in finished statements cannot be modified - financial-accounting closing")

# Insert the blank paragraph that follows it
$blank1 = $newPara.InsertAfter($CR)
$blank1.Font.Size = 24

# Re-apply the run split inside the "After the general ledger..." paragraph so the
# word "quartely" is its own run (matching canonical structure) - same formatting.
$afterStart = $newPara.Start
$prefix = "After the general ledger, the VAT statements and other financial statements were finalized (that can happen monthly, "
$word = "quartely"
$r1 = $tr.Characters($afterStart, $prefix.Length)
$r1.Font.Size = 24
$r2 = $tr.Characters($afterStart + $prefix.Length, $word.Length)
$r2.Font.Size = 24
$r3Start = $afterStart + $prefix.Length + $word.Length
$r3Len = $newPara.Length - $prefix.Length - $word.Length
$r3 = $tr.Characters($r3Start, $r3Len)
$r3.Font.Size = 24

# Fix the body autofit: it should go back to plain normAutofit (no shrink values)
$shp.TextFrame.AutoSize = 1
